$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text so values like "304.75" are not
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.234.79"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.604.91"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "304.75"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "52.38"
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "1.275"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("D11").Value = "0.08159"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "22.94"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "6.607"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "7.385"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "1.607.73"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "94.04"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "18.17"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "6.543"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "12.93"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "23.226.32"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").Value = "2.449"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").Value = "3.076"
$ws.Range("E26").Value = "  +9.66%  "
$ws.Range("D27").Value = "21.19"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "150.06"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "5.284"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "135.36"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "6.742"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "1.778.74"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").Value = "0.9645"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").Value = "0.07496"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").Value = "0.02765"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("D38").Value = "0.2520"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").Value = "0.08802"
$ws.Range("E40").Value = "  -0.32%  "
$ws.Range("D41").Value = "1.410"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("D42").Value = "0.7102"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").Value = "15.72"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").Value = "0.6542"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").Value = "2.333"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "4.010"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "134.10"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").Value = "0.07949"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "1.207"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("E51").Value = "  -2.98%  "
